$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2,I3,I4,I5,I6,I7,I10,I11,I12,I13,I14,I17,I18,I19,I20,I21,I22,I23,I24,I25,I27,I28,I29,I30,I31,I32").NumberFormat = "@"

# Column B (StrID) updates
$ws.Range("B2").Value = 'Bowel_Bag'
$ws.Range("B3").Value = 'SpinalCanal'
$ws.Range("B4").Value = 'BileDuct_Common'
$ws.Range("B5").Value = 'Heart'
$ws.Range("B6").Value = 'CTV'
$ws.Range("B7").Value = 'CTV_High'
$ws.Range("B8").Value = 'CTV_Low'
$ws.Range("B9").Value = 'CTV_Mid'
$ws.Range("B10").Value = 'BODY'
$ws.Range("B11").Value = 'Duodenum'
$ws.Range("B12").Value = 'Esophagus'
$ws.Range("B13").Value = 'Stomach'
$ws.Range("B14").Value = 'GTV'
$ws.Range("B15").Value = 'GTVn'
$ws.Range("B16").Value = 'GTVp'
$ws.Range("B17").Value = 'Liver'
$ws.Range("B18").Value = 'Bowel_Small'
$ws.Range("B19").Value = 'Bowel_Large'
$ws.Range("B20").Value = 'ITV'
$ws.Range("B21").Value = 'SpinalCord'
$ws.Range("B22").Value = 'Chestwall'
$ws.Range("B23").Value = 'PTV'
$ws.Range("B24").Value = 'PTV_High'
$ws.Range("B25").Value = 'PTV_Low'
$ws.Range("B26").Value = 'PTV_Mid'
$ws.Range("B27").Value = 'Lungs'
$ws.Range("B28").Value = 'Lung_L'
$ws.Range("B29").Value = 'Kidneys'
$ws.Range("B30").Value = 'Kidney_R'
$ws.Range("B31").Value = 'Kidney_L'
$ws.Range("B32").Value = 'GallBladder'

# Column C (Type) updates
$ws.Range("C3").Value = 'ORGAN'
$ws.Range("C6").Value = 'CTV'
$ws.Range("C7").Value = 'CTV'
$ws.Range("C10").Value = 'EXTERNAL'
$ws.Range("C11").Value = 'ORGAN'
$ws.Range("C14").Value = 'GTV'
$ws.Range("C17").Value = 'ORGAN'
$ws.Range("C19").Value = 'ORGAN'
$ws.Range("C20").Value = 'CTV'
$ws.Range("C23").Value = 'PTV'
$ws.Range("C24").Value = 'PTV'
$ws.Range("C25").Value = 'PTV'
$ws.Range("C27").Value = 'ORGAN'
$ws.Range("C28").Value = 'ORGAN'
$ws.Range("C29").Value = 'ORGAN'

# Column E (Color) updates
$ws.Range("E2").Value = 'RGB165 42 42'
$ws.Range("E3").Value = 'RGB138255173'
$ws.Range("E4").Value = 'RGB  0240  0'
$ws.Range("E5").Value = 'RGB255  0  0'
$ws.Range("E6").Value = 'RGB255  0  0'
$ws.Range("E7").Value = 'RGB255  0  0'
$ws.Range("E8").Value = 'RGB255255  0'
$ws.Range("E9").Value = 'RGB  6 82255'
$ws.Range("E10").Value = 'RGB  0255  0'
$ws.Range("E11").Value = 'RGB233 67 67'
$ws.Range("E12").Value = 'RGB255165  0'
$ws.Range("E13").Value = 'RGB164  0  0'
$ws.Range("E14").Value = 'RGB  0240  0'
$ws.Range("E15").Value = 'RGB  0240  0'
$ws.Range("E16").Value = 'RGB  0240  0'
$ws.Range("E17").Value = 'RGB255165  0'
$ws.Range("E18").Value = 'RGB164164  0'
$ws.Range("E19").Value = 'RGB  0240  0'
$ws.Range("E20").Value = 'RGB255165  0'
$ws.Range("E21").Value = 'RGB  0240  0'
$ws.Range("E22").Value = 'RGB  0255255'
$ws.Range("E23").Value = 'RGB255  0  0'
$ws.Range("E24").Value = 'RGB255  0  0'
$ws.Range("E25").Value = 'RGB255255  0'
$ws.Range("E26").Value = 'RGB  6 82255'
$ws.Range("E27").Value = 'RGB  6 82255'
$ws.Range("E28").Value = 'RGB127255212'
$ws.Range("E29").Value = 'RGB  0119170'
$ws.Range("E30").Value = 'RGB138255173'
$ws.Range("E31").Value = 'RGB255255  0'
$ws.Range("E32").Value = 'RGB138255173'
$ws.Range("E33").Value = 'RGB  0127255'

# Column I (iCode) updates (kept as text)
$ws.Range("I2").Value = '99'
$ws.Range("I3").Value = '7647'
$ws.Range("I4").Value = '71892'
$ws.Range("I5").Value = '7088'
$ws.Range("I6").Value = '88'
$ws.Range("I7").Value = '88'
$ws.Range("I10").Value = '256135'
$ws.Range("I11").Value = '7206'
$ws.Range("I12").Value = '7131'
$ws.Range("I13").Value = '7148'
$ws.Range("I14").Value = '88'
$ws.Range("I17").Value = '7197'
$ws.Range("I18").Value = '7200'
$ws.Range("I19").Value = '7201'
$ws.Range("I20").Value = '88'
$ws.Range("I21").Value = '7647'
$ws.Range("I22").Value = '50060'
$ws.Range("I23").Value = '88'
$ws.Range("I24").Value = '88'
$ws.Range("I25").Value = '88'
$ws.Range("I27").Value = '7195'
$ws.Range("I28").Value = '7310'
$ws.Range("I29").Value = '7203'
$ws.Range("I30").Value = '7204'
$ws.Range("I31").Value = '7205'
$ws.Range("I32").Value = '7202'

Write-Host "Applied Xaml_Maker test.xlsx updates"
